# Update the "2. Data reporter" section (B6:B10) with the new contact
# details for the National Statistical Committee of the Kyrgyz Republic.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
# Re-apply the font on the organization cell (matches how Excel stamped a
# fresh font/style entry on this cell when the text was retyped).
$ws.Range("B6").Font.Name = "Calibri"

$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# The author's last selection before saving was B10.
$ws.Range("B10").Select()
